# Add a new "2023" column (column R) to the table, mirroring the
# formatting of the existing "2022" column (column Q), and fill in the
# known values for the data rows. Rows that only carry a placeholder
# ("…") in column Q get the same placeholder in column R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-ColumnFormat {
    param($rowNumber)
    $ws.Range("Q$rowNumber").Copy()
    $ws.Range("R$rowNumber").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# Row 3 - thick-bottom-border spacer row directly above the header row.
Copy-ColumnFormat 3

# Row 4 - year header row; Q4 holds 2022, R4 should hold 2023.
Copy-ColumnFormat 4
$ws.Range("R4").Value = 2023

# Row 5 - "Total" row.
Copy-ColumnFormat 5
$ws.Range("R5").Value = 11357

# Row 6 - "By sex:" section header (no value).
Copy-ColumnFormat 6

# Row 7 - "Women"
Copy-ColumnFormat 7
$ws.Range("R7").Value = 11002

# Row 8 - "Men"
Copy-ColumnFormat 8
$ws.Range("R8").Value = 355

# Row 9 - "By age group, years:" section header (no value).
Copy-ColumnFormat 9

# Rows 10-24 - detail rows whose 2020/2021/2022 columns only ever held the
# "…" placeholder string; column R gets the same placeholder.
for ($r = 10; $r -le 24; $r++) {
    Copy-ColumnFormat $r
    $ws.Range("R$r").Value = $ws.Range("Q$r").Value()
}

# Row 25 - "By occupation:" total-style row (bottom border), same placeholder.
Copy-ColumnFormat 25
$ws.Range("R25").Value = $ws.Range("Q25").Value()

# Reset the selection away from the stale Q3 reference left over from
# before the new column existed.
$ws.Range("A1").Select()
